# Update '想去人数' (F column) figures per commit 456a3b4 gh-pages data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 354
$ws.Range("F6").Value = 515
$ws.Range("F9").Value = 259
$ws.Range("F10").Value = 372
$ws.Range("F11").Value = 339
$ws.Range("F12").Value = 581
$ws.Range("F13").Value = 732
$ws.Range("F14").Value = 1486
$ws.Range("F15").Value = 1486
$ws.Range("F18").Value = 1337
$ws.Range("F20").Value = 248
$ws.Range("F21").Value = 238
$ws.Range("F22").Value = 12
$ws.Range("F24").Value = 6434
$ws.Range("F25").Value = 4778
$ws.Range("F26").Value = 130
$ws.Range("F28").Value = 190
$ws.Range("F29").Value = 115
$ws.Range("F32").Value = 1247
$ws.Range("F33").Value = 182
$ws.Range("F34").Value = 35
$ws.Range("F35").Value = 583
$ws.Range("F37").Value = 1333
$ws.Range("F38").Value = 223
$ws.Range("F40").Value = 136
$ws.Range("F41").Value = 56
$ws.Range("F42").Value = 87
$ws.Range("F43").Value = 85

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 37
$ws.Range("F7").Value = 19
$ws.Range("F15").Value = 234

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2427
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 171
$ws.Range("F8").Value = 37
$ws.Range("F9").Value = 354
$ws.Range("F10").Value = 515
$ws.Range("F13").Value = 259
$ws.Range("F15").Value = 372
$ws.Range("F16").Value = 339
$ws.Range("F17").Value = 581
$ws.Range("F18").Value = 732
$ws.Range("F19").Value = 1486
$ws.Range("F20").Value = 1486
$ws.Range("F23").Value = 1337
$ws.Range("F25").Value = 248
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 5
$ws.Range("F29").Value = 37
$ws.Range("F30").Value = 19
$ws.Range("F31").Value = 6434
$ws.Range("F32").Value = 4778
$ws.Range("F33").Value = 130
$ws.Range("F34").Value = 190
$ws.Range("F36").Value = 1247
$ws.Range("F37").Value = 182
$ws.Range("F38").Value = 35
$ws.Range("F40").Value = 583
$ws.Range("F44").Value = 1333
$ws.Range("F45").Value = 223
$ws.Range("F46").Value = 136
$ws.Range("F47").Value = 56
$ws.Range("F48").Value = 87
$ws.Range("F49").Value = 234
